# "Fruta / hortaliza, semanal"
#
# A new weekly price-report row is inserted as row 49 on the single
# worksheet. Every row that used to be 49..105 shifts down by one (to
# 50..106), which is exactly what happens when a real row is inserted
# above the old row 49 in Excel - no other existing cell needs to be
# rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 49; everything below (old rows
# 49-105) shifts down to 50-106, growing the used range to A1:R106.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record.
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44874
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 100112026
$ws.Range("G49").Value = "Haba"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 7000
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = 7000
$ws.Range("N49").Value = '$/saco 25 kilos'
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 280
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
